$wb = $excel.ActiveWorkbook

# Old / new identifiers used throughout the workbook.
$oldGuid = "fbcf006e-ec61-4415-9760-789c0a0c4efa"
$newGuid = "79f4ce9c-39ea-47b5-9073-169b98780371"
$oldHash = "9c674f85aaad76d4c810ef6dbea7a3345e6ff041"
$newHash = "b70fee023c1931b45a75392b25e905039efcd570"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuid + ".md"
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"
$wsOverview.Range("G2").Value = "2016-10-18 12:44:47"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\" + $newGuid + ".md"
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuid + ".md"
$wsZhCn.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-10-18 12:44:36"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = $newGuid + ".md"
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuid + ".md"
$wsDeDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-10-18 12:44:47"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = $newGuid + ".md"
}
